$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = "a2"
$ws.Range("B3").Value = "c2"
$ws.Range("B4").Value = "c2"
$ws.Range("B5").Value = "d2"
$ws.Range("B6").Value = "e2"
$ws.Range("A6").Value = "abcdefg"

$ws.Range("E6").Select()
